$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1756756756756757
$ws.Range("C2").Value = 0.5837837837837838
$ws.Range("J2").Value = 0.01621621621621622
$ws.Range("P2").Value = 0.1324324324324324
$ws.Range("S2").Value = 0.0918918918918919
$ws.Range("B3").Value = 0.01357466063348416
$ws.Range("C3").Value = 0.009049773755656109
$ws.Range("J3").Value = 0.01357466063348416
$ws.Range("P3").Value = 0.7104072398190046
$ws.Range("S3").Value = 0.253393665158371
$ws.Range("J4").Value = 0.02
$ws.Range("O4").Value = 0.02
$ws.Range("P4").Value = 0.72
$ws.Range("S4").Value = 0.24
$ws.Range("B6").Value = 0.06018518518518518
$ws.Range("D6").Value = 0.02314814814814815
$ws.Range("F6").Value = 0.04166666666666666
$ws.Range("J6").Value = 0.2175925925925926
$ws.Range("O6").Value = 0.004629629629629629
$ws.Range("Q6").Value = 0.125
$ws.Range("R6").Value = 0.06944444444444445
$ws.Range("S6").Value = 0.4583333333333333
$ws.Range("B7").Value = 0.1400966183574879
$ws.Range("D7").Value = 0.01932367149758454
$ws.Range("F7").Value = 0.04830917874396135
$ws.Range("J7").Value = 0.08695652173913043
$ws.Range("O7").Value = 0.01932367149758454
$ws.Range("Q7").Value = 0.1352657004830918
$ws.Range("R7").Value = 0.05314009661835749
$ws.Range("S7").Value = 0.4975845410628019
$ws.Range("B8").Value = 0.1290322580645161
$ws.Range("D8").Value = 0.02016129032258064
$ws.Range("E8").Value = 0.004032258064516129
$ws.Range("F8").Value = 0.07056451612903226
$ws.Range("J8").Value = 0.0846774193548387
$ws.Range("O8").Value = 0.01814516129032258
$ws.Range("Q8").Value = 0.1391129032258064
$ws.Range("R8").Value = 0.07459677419354839
$ws.Range("S8").Value = 0.4596774193548387
$ws.Range("B9").Value = 0.1736111111111111
$ws.Range("D9").Value = 0.02777777777777778
$ws.Range("F9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.09027777777777778
$ws.Range("O9").Value = 0.006944444444444444
$ws.Range("Q9").Value = 0.1597222222222222
$ws.Range("R9").Value = 0.04861111111111111
$ws.Range("S9").Value = 0.4375
$ws.Range("B10").Value = 0.1454545454545454
$ws.Range("D10").Value = 0.02510822510822511
$ws.Range("F10").Value = 0.08051948051948052
$ws.Range("J10").Value = 0.103030303030303
$ws.Range("O10").Value = 0.01558441558441558
$ws.Range("Q10").Value = 0.1662337662337662
$ws.Range("R10").Value = 0.07012987012987013
$ws.Range("S10").Value = 0.3939393939393939
$ws.Range("G11").Value = 0.1478260869565217
$ws.Range("J11").Value = 0.1072463768115942
$ws.Range("K11").Value = 0.2202898550724638
$ws.Range("L11").Value = 0.5130434782608696
$ws.Range("S11").Value = 0.01159420289855072
$ws.Range("G12").Value = 0.6968085106382979
$ws.Range("J12").Value = 0.1914893617021277
$ws.Range("K12").Value = 0.01063829787234043
$ws.Range("L12").Value = 0.02659574468085106
$ws.Range("S12").Value = 0.07446808510638298
$ws.Range("G13").Value = 0.7111111111111111
$ws.Range("J13").Value = 0.2666666666666667
$ws.Range("S13").Value = 0.02222222222222222
$ws.Range("F15").Value = 0.03045685279187817
$ws.Range("H15").Value = 0.1319796954314721
$ws.Range("I15").Value = 0.08121827411167512
$ws.Range("J15").Value = 0.3604060913705584
$ws.Range("K15").Value = 0.06598984771573604
$ws.Range("M15").Value = 0.01015228426395939
$ws.Range("O15").Value = 0.05583756345177665
$ws.Range("S15").Value = 0.2639593908629442
$ws.Range("F16").Value = 0.01276595744680851
$ws.Range("H16").Value = 0.1872340425531915
$ws.Range("I16").Value = 0.07234042553191489
$ws.Range("J16").Value = 0.3702127659574468
$ws.Range("K16").Value = 0.1148936170212766
$ws.Range("M16").Value = 0.01702127659574468
$ws.Range("O16").Value = 0.05957446808510639
$ws.Range("S16").Value = 0.1659574468085106
$ws.Range("F17").Value = 0.01474926253687316
$ws.Range("H17").Value = 0.2035398230088496
$ws.Range("I17").Value = 0.06194690265486726
$ws.Range("J17").Value = 0.4070796460176991
$ws.Range("K17").Value = 0.08554572271386431
$ws.Range("M17").Value = 0.02949852507374631
$ws.Range("O17").Value = 0.03244837758112094
$ws.Range("S17").Value = 0.1651917404129793
$ws.Range("F18").Value = 0.02666666666666667
$ws.Range("H18").Value = 0.2133333333333333
$ws.Range("I18").Value = 0.04666666666666667
$ws.Range("J18").Value = 0.4066666666666667
$ws.Range("K18").Value = 0.12
$ws.Range("O18").Value = 0.08
$ws.Range("S18").Value = 0.1066666666666667
$ws.Range("F19").Value = 0.01199717713479181
$ws.Range("H19").Value = 0.2279463655610445
$ws.Range("I19").Value = 0.05928016937191249
$ws.Range("J19").Value = 0.3394495412844037
$ws.Range("K19").Value = 0.1284403669724771
$ws.Range("M19").Value = 0.0218772053634439
$ws.Range("N19").Value = 0.0007057163020465773
$ws.Range("O19").Value = 0.06563161609033169
$ws.Range("S19").Value = 0.1446718419195483
